$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- add the two new sheets in order, right after sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "23JUL24_SWEEP_ON_AIR"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "20JUL24_SWEEP_ON_NITROGEN"

# --- header row, sheet 2 (23JUL24_SWEEP_ON_AIR) ---
$ws2.Cells.Item(1,1).Value = "f_add_raw"
$ws2.Cells.Item(1,2).Value = "f_name_raw"
$ws2.Cells.Item(1,3).Value = "Qsh_DMA (lpm)"
$ws2.Cells.Item(1,4).Value = "DF"
$ws2.Cells.Item(1,5).Value = "sid"

# --- header row, sheet 3 (20JUL24_SWEEP_ON_NITROGEN) ---
$ws3.Cells.Item(1,1).Value = "f_add_raw"
$ws3.Cells.Item(1,2).Value = "f_name_raw"
$ws3.Cells.Item(1,3).Value = "Qsh_DMA (lpm)"
$ws3.Cells.Item(1,4).Value = "DF"
$ws3.Cells.Item(1,5).Value = "sid"

# --- data rows: write B (date sid) & D (DF ratio) columns first for both
# sheets (in this order) so new shared strings land at the indices the
# target workbook expects, then backfill A/C/E afterwards. ---

# sheet2 B/D columns
$ws2.Cells.Item(2,2).Value = "2024-07-23_105221_SMPS"
$ws2.Cells.Item(2,4).Value = "1:3"
$ws2.Cells.Item(3,2).Value = "2024-07-23_105221_SMPS"
$ws2.Cells.Item(3,4).Value = "5:7"
$ws2.Cells.Item(4,2).Value = "2024-07-23_114958_SMPS"
$ws2.Cells.Item(4,4).Value = "1:3"
$ws2.Cells.Item(5,2).Value = "2024-07-23_114958_SMPS"
$ws2.Cells.Item(5,4).Value = "5:7"
$ws2.Cells.Item(6,2).Value = "2024-07-23_143207_SMPS"
$ws2.Cells.Item(6,4).Value = "1:3"
$ws2.Cells.Item(7,2).Value = "2024-07-23_180812_SMPS"
$ws2.Cells.Item(7,4).Value = "1:3"
$ws2.Cells.Item(8,2).Value = "2024-07-23_155349_SMPS"
$ws2.Cells.Item(8,4).Value = "1:3"
$ws2.Cells.Item(9,2).Value = "2024-07-23_193542_SMPS"
$ws2.Cells.Item(9,4).Value = "1:3"
$ws2.Cells.Item(10,2).Value = "2024-07-23_145054_SMPS"
$ws2.Cells.Item(10,4).Value = "1:3"
$ws2.Cells.Item(11,2).Value = "2024-07-23_183014_SMPS"
$ws2.Cells.Item(11,4).Value = "1:3"
$ws2.Cells.Item(12,2).Value = "2024-07-23_161641_SMPS"
$ws2.Cells.Item(12,4).Value = "1:3"
$ws2.Cells.Item(13,2).Value = "2024-07-23_195749_SMPS"
$ws2.Cells.Item(13,4).Value = "1:3"
$ws2.Cells.Item(14,2).Value = "2024-07-23_150820_SMPS"
$ws2.Cells.Item(14,4).Value = "1:3"
$ws2.Cells.Item(15,2).Value = "2024-07-23_185212_SMPS"
$ws2.Cells.Item(15,4).Value = "1:3"
$ws2.Cells.Item(16,2).Value = "2024-07-23_164902_SMPS"
$ws2.Cells.Item(16,4).Value = "1:3"
$ws2.Cells.Item(17,2).Value = "2024-07-23_201626_SMPS"
$ws2.Cells.Item(17,4).Value = "1:3"

# sheet3 B/D columns
$ws3.Cells.Item(2,2).Value = "2024-07-20_210133_SMPS"
$ws3.Cells.Item(2,4).Value = "1:3"
$ws3.Cells.Item(3,2).Value = "2024-07-20_213658_SMPS"
$ws3.Cells.Item(3,4).Value = "1:3"
$ws3.Cells.Item(4,2).Value = "2024-07-20_224040_SMPS"
$ws3.Cells.Item(4,4).Value = "1:3"
$ws3.Cells.Item(5,2).Value = "2024-07-20_222608_SMPS"
$ws3.Cells.Item(5,4).Value = "1:3"
$ws3.Cells.Item(6,2).Value = "2024-07-20_230642_SMPS"
$ws3.Cells.Item(6,4).Value = "1:3"
$ws3.Cells.Item(7,2).Value = "2024-07-20_232313_SMPS"
$ws3.Cells.Item(7,4).Value = "1:3"
$ws3.Cells.Item(8,2).Value = "2024-07-20_234336_SMPS"
$ws3.Cells.Item(8,4).Value = "1:3"

# --- A2 path string (shared between both new sheets; written last so it
# becomes the final new shared-string entry) ---
$ws2.Cells.Item(2,1).Value = "D:\Hamed\CND\PhD\Experiments\PFA-RH122\PFA results\SMPS\SMPS_laptop_datafiles_DEC24\Raw"
$ws3.Cells.Item(2,1).Value = "D:\Hamed\CND\PhD\Experiments\PFA-RH122\PFA results\SMPS\SMPS_laptop_datafiles_DEC24\Raw"

# --- remaining numeric columns (C, E) ---
# sheet2
$ws2.Cells.Item(2,3).Value = 2.5
$ws2.Cells.Item(2,5).Value = 20
$ws2.Cells.Item(3,3).Value = 2.5
$ws2.Cells.Item(3,5).Value = 8
$ws2.Cells.Item(4,3).Value = 2.5
$ws2.Cells.Item(4,5).Value = 20
$ws2.Cells.Item(5,3).Value = 2.5
$ws2.Cells.Item(5,5).Value = 8
$ws2.Cells.Item(6,3).Value = 2.5
$ws2.Cells.Item(6,5).Value = 20
$ws2.Cells.Item(7,3).Value = 2.5
$ws2.Cells.Item(7,5).Value = 8
$ws2.Cells.Item(8,3).Value = 4.5
$ws2.Cells.Item(8,5).Value = 20
$ws2.Cells.Item(9,3).Value = 2.5
$ws2.Cells.Item(9,5).Value = 8
$ws2.Cells.Item(10,3).Value = 2.5
$ws2.Cells.Item(10,5).Value = 20
$ws2.Cells.Item(11,3).Value = 2.5
$ws2.Cells.Item(11,5).Value = 8
$ws2.Cells.Item(12,3).Value = 3
$ws2.Cells.Item(12,5).Value = 20
$ws2.Cells.Item(13,3).Value = 2.5
$ws2.Cells.Item(13,5).Value = 8
$ws2.Cells.Item(14,3).Value = 2.5
$ws2.Cells.Item(14,5).Value = 20
$ws2.Cells.Item(15,3).Value = 2.5
$ws2.Cells.Item(15,5).Value = 8
$ws2.Cells.Item(16,3).Value = 4.5
$ws2.Cells.Item(16,5).Value = 20
$ws2.Cells.Item(17,3).Value = 2.5
$ws2.Cells.Item(17,5).Value = 8

# sheet3
$ws3.Cells.Item(2,3).Value = 2.5
$ws3.Cells.Item(2,5).Value = 20
$ws3.Cells.Item(3,3).Value = 2.5
$ws3.Cells.Item(3,5).Value = 20
$ws3.Cells.Item(4,3).Value = 2.5
$ws3.Cells.Item(4,5).Value = 20
$ws3.Cells.Item(5,3).Value = 2.5
$ws3.Cells.Item(5,5).Value = 20
$ws3.Cells.Item(6,3).Value = 2.5
$ws3.Cells.Item(6,5).Value = 20
$ws3.Cells.Item(7,3).Value = 2.5
$ws3.Cells.Item(7,5).Value = 20
$ws3.Cells.Item(8,3).Value = 2.5
$ws3.Cells.Item(8,5).Value = 20

# --- selections / active-sheet state to mirror the target workbook ---
$ws3.Range("B2:B8").Select()
$ws2.Range("B2:B17").Select()
$ws2.Activate()

Write-Host "done"
